# Update gh-pages to output generated at 456a3b4
# Applies updated "want-to-go" (F) and "min-price" (G) figures across the
# four sheets of the workbook.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions)
$wsExhibit.Range("F2").Value  = 222
$wsExhibit.Range("F4").Value  = 403
$wsExhibit.Range("G6").Value  = 75
$wsExhibit.Range("F7").Value  = 113
$wsExhibit.Range("F8").Value  = 10312
$wsExhibit.Range("G8").Value  = 95
$wsExhibit.Range("F10").Value = 3562
$wsExhibit.Range("F14").Value = 2847
$wsExhibit.Range("F17").Value = 2191
$wsExhibit.Range("F20").Value = 31
$wsExhibit.Range("F21").Value = 396
$wsExhibit.Range("F23").Value = 158
$wsExhibit.Range("F26").Value = 239
$wsExhibit.Range("F28").Value = 1329
$wsExhibit.Range("F29").Value = 17
$wsExhibit.Range("F30").Value = 1262
$wsExhibit.Range("F34").Value = 3871
$wsExhibit.Range("F35").Value = 3275
$wsExhibit.Range("F38").Value = 1048
$wsExhibit.Range("F39").Value = 407
$wsExhibit.Range("F42").Value = 111
$wsExhibit.Range("F44").Value = 75
$wsExhibit.Range("F47").Value = 15

# 演出 (Shows)
$wsShow.Range("F7").Value  = 3
$wsShow.Range("F16").Value = 182

# 本地生活 (Local life)
$wsLocal.Range("F3").Value = 994
$wsLocal.Range("F4").Value = 131
$wsLocal.Range("F5").Value = 2078

# 全部类型 (All types)
$wsAll.Range("F3").Value  = 994
$wsAll.Range("F4").Value  = 131
$wsAll.Range("F6").Value  = 403
$wsAll.Range("G9").Value  = 75
$wsAll.Range("F10").Value = 113
$wsAll.Range("F11").Value = 10312
$wsAll.Range("G11").Value = 95
$wsAll.Range("F13").Value = 3562
$wsAll.Range("F19").Value = 2191
$wsAll.Range("F22").Value = 31
$wsAll.Range("F23").Value = 396
$wsAll.Range("F24").Value = 158
$wsAll.Range("F26").Value = 239
$wsAll.Range("F27").Value = 1329
$wsAll.Range("F28").Value = 17
$wsAll.Range("F29").Value = 1262
$wsAll.Range("F32").Value = 3
$wsAll.Range("F36").Value = 3275
$wsAll.Range("F37").Value = 1048
$wsAll.Range("F45").Value = 111
$wsAll.Range("F46").Value = 75
$wsAll.Range("F48").Value = 15
$wsAll.Range("F49").Value = 182
